# Apply updated NATMI Gnai2-Adra2a LR-pair values (Dr Hou advice revision)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: E F G H I J K L M N O P Q R S T (numeric data columns, A-D are text/labels and unchanged)
# Row 2 (ECs -> ECs)
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 184.0626906666667
$ws.Range("H2").Value = 552.188072
$ws.Range("I2").Value = 0.6510505751503485
$ws.Range("J2").Value = 0.6510505751503486
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.6723143333333333
$ws.Range("N2").Value = 2.016943
$ws.Range("O2").Value = 0.1884809556899308
$ws.Range("P2").Value = 0.1884809556899308
$ws.Range("Q2").Value = 123.7479851670996
$ws.Range("R2").Value = 1113.731866503896
$ws.Range("S2").Value = 0.1227106346068168
$ws.Range("T2").Value = 0.1227106346068168

# Row 3 (ECs -> FAPs)
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 184.0626906666667
$ws.Range("H3").Value = 552.188072
$ws.Range("I3").Value = 0.6510505751503485
$ws.Range("J3").Value = 0.6510505751503486
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.945804
$ws.Range("N3").Value = 2.837412
$ws.Range("O3").Value = 0.2651528206032981
$ws.Range("P3").Value = 0.2651528206032981
$ws.Range("Q3").Value = 174.087229083296
$ws.Range("R3").Value = 1566.785061749664
$ws.Range("S3").Value = 0.1726278963565144
$ws.Range("T3").Value = 0.1726278963565145

# Row 4 (ECs -> sCs)
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 184.0626906666667
$ws.Range("H4").Value = 552.188072
$ws.Range("I4").Value = 0.6510505751503485
$ws.Range("J4").Value = 0.6510505751503486
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.948896333333333
$ws.Range("N4").Value = 5.846689
$ws.Range("O4").Value = 0.546366223706771
$ws.Range("P4").Value = 0.546366223706771
$ws.Range("Q4").Value = 358.7191029437342
$ws.Range("R4").Value = 3228.471926493608
$ws.Range("S4").Value = 0.3557120441870172
$ws.Range("T4").Value = 0.3557120441870173

# Row 5 (FAPs -> ECs)
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 57.4434
$ws.Range("H5").Value = 172.3302
$ws.Range("I5").Value = 0.2031838091312023
$ws.Range("J5").Value = 0.2031838091312023
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.6723143333333333
$ws.Range("N5").Value = 2.016943
$ws.Range("O5").Value = 0.1884809556899308
$ws.Range("P5").Value = 0.1884809556899308
$ws.Range("Q5").Value = 38.6200211754
$ws.Range("R5").Value = 347.5801905786
$ws.Range("S5").Value = 0.03829627852576951
$ws.Range("T5").Value = 0.0382962785257695

# Row 6 (FAPs -> FAPs)
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 57.4434
$ws.Range("H6").Value = 172.3302
$ws.Range("I6").Value = 0.2031838091312023
$ws.Range("J6").Value = 0.2031838091312023
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.945804
$ws.Range("N6").Value = 2.837412
$ws.Range("O6").Value = 0.2651528206032981
$ws.Range("P6").Value = 0.2651528206032981
$ws.Range("Q6").Value = 54.3301974936
$ws.Range("R6").Value = 488.9717774424
$ws.Range("S6").Value = 0.05387476009206045
$ws.Range("T6").Value = 0.05387476009206045

# Row 7 (FAPs -> sCs)
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 57.4434
$ws.Range("H7").Value = 172.3302
$ws.Range("I7").Value = 0.2031838091312023
$ws.Range("J7").Value = 0.2031838091312023
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.948896333333333
$ws.Range("N7").Value = 5.846689
$ws.Range("O7").Value = 0.546366223706771
$ws.Range("P7").Value = 0.546366223706771
$ws.Range("Q7").Value = 111.9512316342
$ws.Range("R7").Value = 1007.5610847078
$ws.Range("S7").Value = 0.1110127705133723
$ws.Range("T7").Value = 0.1110127705133723

# Row 8 (sCs -> ECs)
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 41.21033366666666
$ws.Range("H8").Value = 123.631001
$ws.Range("I8").Value = 0.1457656157184491
$ws.Range("J8").Value = 0.1457656157184491
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.6723143333333333
$ws.Range("N8").Value = 2.016943
$ws.Range("O8").Value = 0.1884809556899308
$ws.Range("P8").Value = 0.1884809556899308
$ws.Range("Q8").Value = 27.70629800554922
$ws.Range("R8").Value = 249.356682049943
$ws.Range("S8").Value = 0.02747404255734449
$ws.Range("T8").Value = 0.02747404255734449

# Row 9 (sCs -> FAPs)
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 41.21033366666666
$ws.Range("H9").Value = 123.631001
$ws.Range("I9").Value = 0.1457656157184491
$ws.Range("J9").Value = 0.1457656157184491
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.945804
$ws.Range("N9").Value = 2.837412
$ws.Range("O9").Value = 0.2651528206032981
$ws.Range("P9").Value = 0.2651528206032981
$ws.Range("Q9").Value = 38.976898423268
$ws.Range("R9").Value = 350.792085809412
$ws.Range("S9").Value = 0.03865016415472323
$ws.Range("T9").Value = 0.03865016415472323

# Row 10 (sCs -> sCs)
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 41.21033366666666
$ws.Range("H10").Value = 123.631001
$ws.Range("I10").Value = 0.1457656157184491
$ws.Range("J10").Value = 0.1457656157184491
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.948896333333333
$ws.Range("N10").Value = 5.846689
$ws.Range("O10").Value = 0.546366223706771
$ws.Range("P10").Value = 0.546366223706771
$ws.Range("Q10").Value = 80.31466817840987
$ws.Range("R10").Value = 722.832013605689
$ws.Range("S10").Value = 0.07964140900638139
$ws.Range("T10").Value = 0.07964140900638139
